$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts the old row 7 down to row 8),
# and copy the current row 6 contents into it.
$ws.Rows.Item(7).Insert()

for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(7, $col).Value2 = $ws.Cells.Item(6, $col).Value2
}
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat

# Update row 6 (original row) with the new date / volume values.
$ws.Cells.Item(6, 4).Value2 = 44757
$ws.Cells.Item(6, 10).Value2 = 30
